$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new daily data row (row 33)
$ws.Range("A33").Value = 92
$ws.Range("B33").Value = 59457
$ws.Range("C33").Value = 8251
$ws.Range("D33").Value = 43
$ws.Range("E33").Value = 187

# Extend predictions / update view state
$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Range("D38").Select()
